$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 422726.6
$ws.Range("I15").Value = 422726.6
$ws.Range("K15").Value = 1268179.8
$ws.Range("M15").Value = -1268010.8
$ws.Range("H64").Value = 19199.6
$ws.Range("I64").Value = 15332.667
$ws.Range("K64").Value = 15332.667
$ws.Range("M64").Value = -15084.667
$ws.Range("H67").Value = 19199.6
$ws.Range("I67").Value = 15332.667
$ws.Range("K67").Value = 15332.667
$ws.Range("M67").Value = -14474.667
$ws.Range("H74").Value = 6694.5
$ws.Range("I74").Value = 6036.5454
$ws.Range("K74").Value = 6036.5454
$ws.Range("M74").Value = -5100.5454
$ws.Range("H77").Value = 6694.5
$ws.Range("I77").Value = 6036.5454
$ws.Range("K77").Value = 30182.727
$ws.Range("M77").Value = -25502.727
$ws.Range("H138").Value = 3336.7817
$ws.Range("I138").Value = 836.3570999999999
$ws.Range("J138").Value = 4190.5854
$ws.Range("K138").Value = 2509.0713
$ws.Range("L138").Value = 12571.7562
$ws.Range("M138").Value = 2630.9287
$ws.Range("N138").Value = -22851.7562
$ws.Range("H141").Value = 2309.111
$ws.Range("I141").Value = 1974
$ws.Range("J141").Value = 4990
$ws.Range("K141").Value = 5922
$ws.Range("L141").Value = 14970
$ws.Range("M141").Value = -742
$ws.Range("N141").Value = -25330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1483.1833
$ws.Range("I32").Value = 551.875
$ws.Range("K32").Value = 551.875
$ws.Range("M32").Value = -264.875
$ws.Range("H45").Value = 1462.5
$ws.Range("I45").Value = 1383.3334
$ws.Range("J45").Value = 1700
$ws.Range("K45").Value = 1383.3334
$ws.Range("L45").Value = 1700
$ws.Range("M45").Value = -1006.3334
$ws.Range("N45").Value = -2454
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H122").Value = 4980.3438
$ws.Range("I122").Value = 3911.375
$ws.Range("K122").Value = 11734.125
$ws.Range("M122").Value = -9284.125
$ws.Range("H132").Value = 33730.105
$ws.Range("I132").Value = 2127.25
$ws.Range("J132").Value = 223347.25
$ws.Range("K132").Value = 6381.75
$ws.Range("L132").Value = 670041.75
$ws.Range("M132").Value = -3851.75
$ws.Range("N132").Value = -675101.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2129
$ws.Range("I11").Value = 525
$ws.Range("K11").Value = 525
$ws.Range("M11").Value = -385
$ws.Range("H86").Value = 4284.769
$ws.Range("I86").Value = 2226.8
$ws.Range("J86").Value = 5571
$ws.Range("K86").Value = 2226.8
$ws.Range("L86").Value = 5571
$ws.Range("M86").Value = -1103.8
$ws.Range("N86").Value = -7817
$ws.Range("H89").Value = 4284.769
$ws.Range("I89").Value = 2226.8
$ws.Range("J89").Value = 5571
$ws.Range("K89").Value = 11134
$ws.Range("L89").Value = 27855
$ws.Range("M89").Value = -5518
$ws.Range("N89").Value = -39087
$ws.Range("H94").Value = 48037.566
$ws.Range("I94").Value = 812.9167
$ws.Range("K94").Value = 812.9167
$ws.Range("M94").Value = -361.9167
$ws.Range("H108").Value = 35000
$ws.Range("J108").Value = 35000
$ws.Range("L108").Value = 35000
$ws.Range("N108").Value = -42680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4105.9165
$ws.Range("I31").Value = 1194
$ws.Range("J31").Value = 6711.316
$ws.Range("K31").Value = 1194
$ws.Range("L31").Value = 6711.316
$ws.Range("M31").Value = -899
$ws.Range("N31").Value = -7301.316
$ws.Range("H34").Value = 4105.9165
$ws.Range("I34").Value = 1194
$ws.Range("J34").Value = 6711.316
$ws.Range("K34").Value = 1194
$ws.Range("L34").Value = 6711.316
$ws.Range("M34").Value = -992
$ws.Range("N34").Value = -7115.316
$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("K35").Value = 1000
$ws.Range("M35").Value = -706
$ws.Range("H62").Value = 8315.223
$ws.Range("I62").Value = 3637.8
$ws.Range("K62").Value = 3637.8
$ws.Range("M62").Value = -3013.8
$ws.Range("H65").Value = 8315.223
$ws.Range("I65").Value = 3637.8
$ws.Range("K65").Value = 18189
$ws.Range("M65").Value = -15069
$ws.Range("H94").Value = 6377.615
$ws.Range("I94").Value = 5063.75
$ws.Range("J94").Value = 6961.5557
$ws.Range("K94").Value = 5063.75
$ws.Range("L94").Value = 6961.5557
$ws.Range("M94").Value = -4612.75
$ws.Range("N94").Value = -7863.5557
$ws.Range("H107").Value = 947.375
$ws.Range("I107").Value = 716.4545000000001
$ws.Range("K107").Value = 716.4545000000001
$ws.Range("M107").Value = 1203.5455
$ws.Range("H122").Value = 4031.4
$ws.Range("I122").Value = 2945
$ws.Range("K122").Value = 8835
$ws.Range("M122").Value = -6385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 116547.05
$ws.Range("I122").Value = 553.1667
$ws.Range("J122").Value = 134861.88
$ws.Range("K122").Value = 4978.5003
$ws.Range("L122").Value = 1213756.92
$ws.Range("M122").Value = -2528.5003
$ws.Range("N122").Value = -1218656.92
$ws.Range("H131").Value = 36113156
$ws.Range("I131").Value = 66667604
$ws.Range("J131").Value = 20835932
$ws.Range("K131").Value = 200002812
$ws.Range("L131").Value = 62507796
$ws.Range("M131").Value = -199997772
$ws.Range("N131").Value = -62517876
$ws.Range("H132").Value = 2662.5454
$ws.Range("I132").Value = 1037.25
$ws.Range("J132").Value = 6996.6665
$ws.Range("K132").Value = 9335.25
$ws.Range("L132").Value = 62969.9985
$ws.Range("M132").Value = -6805.25
$ws.Range("N132").Value = -68029.9985
$ws.Range("H134").Value = 4440.96
$ws.Range("I134").Value = 4440.96
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13322.88
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8252.880000000001
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 40000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 40000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 40000
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = -40560
$ws.Range("H44").Value = 31000
$ws.Range("I44").Value = 31000
$ws.Range("K44").Value = 31000
$ws.Range("M44").Value = -30404
$ws.Range("H50").Value = 40000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 40000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 40000
$ws.Range("M50").Value = $null
$ws.Range("N50").Value = -40996
$ws.Range("H70").Value = 10182.066
$ws.Range("I70").Value = 9388.666999999999
$ws.Range("J70").Value = 10711
$ws.Range("K70").Value = 9388.666999999999
$ws.Range("L70").Value = 10711
$ws.Range("M70").Value = -9118.666999999999
$ws.Range("N70").Value = -11251
$ws.Range("H73").Value = 10182.066
$ws.Range("I73").Value = 9388.666999999999
$ws.Range("J73").Value = 10711
$ws.Range("K73").Value = 9388.666999999999
$ws.Range("L73").Value = 10711
$ws.Range("M73").Value = -8452.666999999999
$ws.Range("N73").Value = -12583
$ws.Range("H122").Value = 14336.714
$ws.Range("I122").Value = 13434.5
$ws.Range("J122").Value = 19750
$ws.Range("K122").Value = 40303.5
$ws.Range("L122").Value = 59250
$ws.Range("M122").Value = -37853.5
$ws.Range("N122").Value = -64150
$ws.Range("H132").Value = 5857.5713
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5857.5713
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17572.7139
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -22632.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1313.7333
$ws.Range("I46").Value = 977.7143
$ws.Range("J46").Value = 1607.75
$ws.Range("K46").Value = 977.7143
$ws.Range("L46").Value = 1607.75
$ws.Range("M46").Value = -789.7143
$ws.Range("N46").Value = -1983.75
$ws.Range("H53").Value = 30000
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null
$ws.Range("H68").Value = 1923.8334
$ws.Range("I68").Value = 1412.375
$ws.Range("K68").Value = 1412.375
$ws.Range("M68").Value = -663.375
$ws.Range("H71").Value = 1923.8334
$ws.Range("I71").Value = 1412.375
$ws.Range("K71").Value = 7061.875
$ws.Range("M71").Value = -3317.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1294.8572
$ws.Range("I81").Value = 710.6667
$ws.Range("J81").Value = 4800
$ws.Range("K81").Value = 1421.3334
$ws.Range("L81").Value = 9600
$ws.Range("M81").Value = -360.3334
$ws.Range("N81").Value = -11722
$ws.Range("H84").Value = 1294.8572
$ws.Range("I84").Value = 710.6667
$ws.Range("J84").Value = 4800
$ws.Range("K84").Value = 7106.666999999999
$ws.Range("L84").Value = 48000
$ws.Range("M84").Value = -1802.666999999999
$ws.Range("N84").Value = -58608
$ws.Range("H122").Value = 2907.6128
$ws.Range("I122").Value = 2801.3845
$ws.Range("K122").Value = 8404.1535
$ws.Range("M122").Value = -5954.1535
